$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the three oldest years of data (2007, 2008, 2009) which occupied rows 2-4.
# This shifts the remaining years (2010..2020) up to rows 2..12.
$ws.Range("A2:A4").EntireRow.Delete() | Out-Null

# Append the new 2021 data as the new last row (row 13).
# Copy the formatting from the cell above (A12) so the new year label matches
# the style used by every other year label in column A.
$ws.Range("A12").Copy() | Out-Null
$ws.Range("A13").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

$ws.Range("A13").Value = "2021年"
$ws.Range("B13").Value = 6053.397781593
$ws.Range("C13").Value = 329487.335947386
$ws.Range("D13").Value = 96164.80218100001
$ws.Range("E13").Value = 9379.9987956975
$ws.Range("F13").Value = 15061.8801150516
$ws.Range("G13").Value = 22588.0511052504
$ws.Range("I13").Value = 950
$ws.Range("J13").Value = 998.2070262442001
$ws.Range("K13").Value = 225102.820995842
$ws.Range("L13").Value = 212867.195240294
$ws.Range("M13").Value = 212392.894113518
$ws.Range("N13").Value = 4125.2220146041
$ws.Range("O13").Value = 128645.470568562
$ws.Range("P13").Value = 15240.6837398648
$ws.Range("R13").Value = 395702.248424123
$ws.Range("S13").Value = 395702.248424123
$ws.Range("T13").Value = 42931.6755272159
$ws.Range("U13").Value = 219.7520266321
$ws.Range("V13").Value = 2855.6269598502
